# Updates to eps-us commit #54625f6
#
# The "thermochemical water splitting" hydrogen production pathway is
# renamed to "hydrocarbon partial oxidation" on the RHPF sheet (used both
# as a column header and as a row label in the fraction matrix).

$wb = $excel.ActiveWorkbook

$oldText = "thermochemical water splitting"
$newText = "hydrocarbon partial oxidation"

foreach ($ws in $wb.Worksheets) {
    [void]$ws.Cells.Replace($oldText, $newText)
}
